$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new daily record at row 345 (shifts existing rows 345.. down by one)
$ws.Rows.Item(345).Insert()
$ws.Range("A345").Value = 5
$ws.Range("B345").Value = "Macroferia Regional de Talca"
$ws.Range("C345").Value = "Maule"
$ws.Range("D345").Value = 45120
$ws.Range("E345").Value = 7
$ws.Range("F345").Value = 100112008
$ws.Range("G345").Value = "Coliflor"
$ws.Range("H345").Value = "Sin especificar"
$ws.Range("I345").Value = "Primera"
$ws.Range("J345").Value = 5000
$ws.Range("K345").Value = 600
$ws.Range("L345").Value = 600
$ws.Range("M345").Value = 600
$ws.Range("N345").Value = '$/unidad'
$ws.Range("O345").Value = "Región del Maule"
$ws.Range("P345").Value = 600
$ws.Range("Q345").Value = 1
$ws.Range("R345").Value = "Hortaliza"

# Insert second new daily record at row 455 (shifts rows 455.. down by one more)
$ws.Rows.Item(455).Insert()
$ws.Range("A455").Value = 5
$ws.Range("B455").Value = "Macroferia Regional de Talca"
$ws.Range("C455").Value = "Maule"
$ws.Range("D455").Value = 45121
$ws.Range("E455").Value = 7
$ws.Range("F455").Value = 100112008
$ws.Range("G455").Value = "Coliflor"
$ws.Range("H455").Value = "Sin especificar"
$ws.Range("I455").Value = "Primera"
$ws.Range("J455").Value = 5000
$ws.Range("K455").Value = 600
$ws.Range("L455").Value = 600
$ws.Range("M455").Value = 600
$ws.Range("N455").Value = '$/unidad'
$ws.Range("O455").Value = "Región del Maule"
$ws.Range("P455").Value = 600
$ws.Range("Q455").Value = 1
$ws.Range("R455").Value = "Hortaliza"
